$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added 64 bit designs" -- refresh the utilization numbers in row 2 with the
# new 64-bit design measurements.
$ws.Range("B2").Value = 50.746238708496094
$ws.Range("C2").Value = 5.344827651977539
$ws.Range("D2").Value = 19.862781524658203
$ws.Range("E2").Value = 57.85714340209961
$ws.Range("F2").Value = 72.7272720336914

# Widen column F (DSP) to match the other numeric columns now that it holds
# wider values.
$ws.Columns.Item(6).ColumnWidth = 10.8
